# Fruta / hortaliza, semanal
# Insert a new weekly record at row 201 ("Zapallo italiano" / Macroferia
# Regional de Talca), shifting the existing rows 201-247 down to 202-248.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 201 - this pushes old rows 201..247
# down to 202..248, preserving their data and formatting.
$ws.Rows("201:201").Insert()

# Populate the new row 201 with the new observation's data.
$ws.Range("A201").Value = 5
$ws.Range("B201").Value = "Macroferia Regional de Talca"
$ws.Range("C201").Value = "Maule"
$ws.Range("D201").Value = 44543
$ws.Range("E201").Value = 7
$ws.Range("F201").Value = 100112032
$ws.Range("G201").Value = "Zapallo italiano"
$ws.Range("H201").Value = "Sin especificar"
$ws.Range("I201").Value = "Primera"
$ws.Range("J201").Value = 500
$ws.Range("K201").Value = 5000
$ws.Range("L201").Value = 5000
$ws.Range("M201").Value = 5000
$ws.Range("N201").Value = "$/caja 60 unidades"
$ws.Range("O201").Value = "Región del Maule"
$ws.Range("P201").Value = 83
$ws.Range("Q201").Value = 60
$ws.Range("R201").Value = "Hortaliza"
